$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")

# Row 20 = columnsEducationE1a: update its value from "20" to "19".
# Leading apostrophe forces a text (quote-prefixed) entry, matching the
# original cell's text storage (style carries the quote-prefix flag) rather
# than letting Excel coerce "19" into a numeric cell.
$ws.Cells.Item(20, 2).Value = "'19"
